# Apply "Add files via upload" edit:
# Remove the "example_3/" sub-section (3 rows collapsed to 1: keep only the
# "dataset for Example 3" row) and remove the "main_example_3.m" row from
# the "copula estimation/" section. Then renumber the running index in
# column A, fix the column B width, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 held the "example_3/" second-level-folder marker (with
# marginal_distribution.R); row 11 held the "main.m" line for example_3.
# Both get removed, leaving only the former row 12 (the dataset row),
# which slides up to become the new row 10.
$ws.Rows("10:11").Delete()

# The row that now sits at position 14 is "main_example_3.m"
# (estimate copula models for Example 3) - remove it too.
$ws.Rows("14").Delete()

# Renumber the sequential index in column A (rows 2..29 => 1..28)
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Adjust column B width (closest achievable value to 25.6640625)
$ws.Columns("B").ColumnWidth = 25

# Update the active cell selection
$ws.Range("C15").Select()
